$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Symbol_table")
$ws.Activate()

$ws.Range("B4").Value = "Q       8.1"
$ws.Range("B2").Value = "I       4.0"
$ws.Range("B3").Value = "I       4.1"

$ws.Range("D8").Select()
